$wb = $excel.ActiveWorkbook

$wsGG = $wb.Worksheets.Item("GG Map API")

# Add the two new worksheets right after "GG Map API", in final order.
$wsOSM = $wb.Worksheets.Add($null, $wsGG)
$wsOSM.Name = "OpenStreetMap"

$wsQGIS = $wb.Worksheets.Add($null, $wsOSM)
$wsQGIS.Name = "QGIS 3"

# --- Populate "OpenStreetMap" sheet (B2, B4, B3 order to match shared-string
# insertion order of the target workbook) ---
$wsOSM.Range("B2").Value = "https://stackoverflow.com/questions/925164/openstreetmap-embedding-map-in-webpage-like-google-maps"
$wsOSM.Hyperlinks.Add($wsOSM.Range("B2"), "https://stackoverflow.com/questions/925164/openstreetmap-embedding-map-in-webpage-like-google-maps") | Out-Null
$wsOSM.Range("B2").Style = "Hyperlink"

# --- Populate "QGIS 3" sheet ---
$wsQGIS.Range("B2").Value = "https://qgis.org/en/docs/index.html"
$wsQGIS.Hyperlinks.Add($wsQGIS.Range("B2"), "https://qgis.org/en/docs/index.html") | Out-Null
$wsQGIS.Range("B2").Style = "Hyperlink"

$wsOSM.Range("B4").Value = "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwjzorWent35AhVflFYBHWpMBdUQFnoECAwQAQ&url=https%3A%2F%2Ftowardsdatascience.com%2Floading-data-from-openstreetmap-with-python-and-the-overpass-api-513882a27fd0&usg=AOvVaw3YLvFD8iUT1pcjWWphxab4"
$wsOSM.Hyperlinks.Add($wsOSM.Range("B4"), "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwjzorWent35AhVflFYBHWpMBdUQFnoECAwQAQ&url=https%3A%2F%2Ftowardsdatascience.com%2Floading-data-from-openstreetmap-with-python-and-the-overpass-api-513882a27fd0&usg=AOvVaw3YLvFD8iUT1pcjWWphxab4") | Out-Null
$wsOSM.Range("B4").Style = "Hyperlink"

$wsOSM.Range("B3").Value = "https://pygis.io/docs/d_access_osm.html"
$wsOSM.Hyperlinks.Add($wsOSM.Range("B3"), "https://pygis.io/docs/d_access_osm.html") | Out-Null
$wsOSM.Range("B3").Style = "Hyperlink"

# "QGIS 3" sheet keeps the cursor on B2 (its only populated cell).
$wsQGIS.Range("B2").Select() | Out-Null

# Active selection: OpenStreetMap is the selected tab, cursor on B3.
$wsOSM.Activate()
$wsOSM.Range("B3").Select() | Out-Null
